$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Targets")

# "Baton Rouge Red" -> "Red" for the sales-center column (C2:C22)
$ws.Range("C2:C22").Value = "Red"

# "Convenience" -> "CR&LT" for the channel column (E2:E43)
$ws.Range("E2:E43").Value = "CR&LT"

# Reflect the author's final on-sheet selection: C17:C22 on the Targets sheet
$ws.Activate()
$ws.Range("C17:C22").Select()
